# Added MX-BBX and MX-DPBX accessories in panel accessories sheet for
# Austria, Italy, Slovakia, Netherlands and Denmark market.
#
# Each country sheet has a list of "Accessories" starting a few rows below
# the header block. We insert two new rows - "MX-DPBX" then "MX-BBX" -
# right after the existing "PR1D2-Unmonitored" row (or, where that row
# doesn't exist yet, right before the trailing "0"/"FB800" + "Accessories"
# rows), matching the rest of the already-populated country sheets.

$wb = $excel.ActiveWorkbook

function Add-Accessories {
    param(
        [string]$SheetName,
        [int]$InsertAtRow
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Activate()
    $ws.Rows.Item($InsertAtRow).Insert()
    $ws.Rows.Item($InsertAtRow).Insert()

    $firstCell = $ws.Cells.Item($InsertAtRow, 1)
    $secondCell = $ws.Cells.Item($InsertAtRow + 1, 1)

    $firstCell.Value = "MX-DPBX"
    $secondCell.Value = "MX-BBX"

    $newRange = $ws.Range($firstCell, $secondCell)
    # Match the thin-border style ("s=3") used by the sibling accessory rows.
    $newRange.Borders.ColorIndex = 1
    $newRange.Borders.LineStyle = 1

    $newRange.Select()
}

# Netherlands & Slovakia already have a "PR1D2-Unmonitored" (row 10) entry,
# so the new rows land between it and the following rows.
Add-Accessories "Slovakia" 10
Add-Accessories "Netherlands" 10

# Italy also has that row, but it currently sits one row lower (row 10 holds
# a different accessory), so the new rows are inserted at row 11.
Add-Accessories "Italy" 11

# Austria & Denmark don't have that extra row, so the new entries are
# inserted right at row 10.
Add-Accessories "Austria" 10
Add-Accessories "Denmark" 10

# Update the Spain selection to match the saved state, then make Poland the
# active sheet/tab (this also flips tabSelected from Spain to Poland and
# updates the workbook's activeTab).
$wsSpain = $wb.Worksheets.Item("Spain")
$wsSpain.Range("A11:A12").Select()

$wsPoland = $wb.Worksheets.Item("Poland")
$wsPoland.Activate()
